# Update performance document: add a new "v1243" test-run column (F on
# PartOfSponza, C on Sponza/ComplexMesh) with measured results, and move
# the active/selected sheet & cell around per the new review pass.

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("PartOfSponza")
$sheet2 = $wb.Worksheets.Item("Sponza")
$sheet3 = $wb.Worksheets.Item("ComplexMesh")

# --- PartOfSponza: new column F = "v1243" ---
$sheet1.Range("F1").Value = "v1243"
$sheet1.Range("F2").Value = 102
$sheet1.Range("F3").Value = 102
$sheet1.Range("F4").Value = 103
$sheet1.Range("F5").Value = 101
$sheet1.Range("F6").Value = 103
$sheet1.Range("F7").Value = 101
$sheet1.Range("F8").Value = 102
$sheet1.Range("F9").Value = 102
$sheet1.Range("F10").Value = 102
$sheet1.Range("F11").Value = 101

# --- Sponza: new column C = "v1243" ---
$sheet2.Range("C1").Value = "v1243"
$sheet2.Range("C2").Value = 8844
$sheet2.Range("C3").Value = 8917
$sheet2.Range("C4").Value = 8703
$sheet2.Range("C5").Value = 9020
$sheet2.Range("C6").Value = 9197
$sheet2.Range("C7").Value = 8943
$sheet2.Range("C8").Value = 9048
$sheet2.Range("C9").Value = 9011
$sheet2.Range("C10").Value = 9055
$sheet2.Range("C11").Value = 8898

# --- ComplexMesh: new column C = "v1243" ---
$sheet3.Range("C1").Value = "v1243"
$sheet3.Range("C2").Value = 5553
$sheet3.Range("C3").Value = 5561
$sheet3.Range("C4").Value = 5600
$sheet3.Range("C5").Value = 5607
$sheet3.Range("C6").Value = 5577
$sheet3.Range("C7").Value = 5598
$sheet3.Range("C8").Value = 5617
$sheet3.Range("C9").Value = 5560
$sheet3.Range("C10").Value = 5563
$sheet3.Range("C11").Value = 5573

# --- Selections / active sheet ---
$sheet1.Range("F12").Select() | Out-Null
$sheet2.Range("C14").Select() | Out-Null

$sheet3.Range("D9").Select() | Out-Null
$sheet3.Activate() | Out-Null
